# Auto-generated Excel COM-interop script applying the Universalis price refresh
# described by the commit "chore: update Sheets via scheduled runner".
# Each block updates the currentAveragePrice* / Leve price / profit columns (H:N)
# for the affected leve rows across the ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3460.762
$ws.Range("I74").Value = 971.0345
$ws.Range("J74").Value = 9014.77
$ws.Range("K74").Value = 971.0345
$ws.Range("L74").Value = 9014.77
$ws.Range("M74").Value = -97.03449999999998
$ws.Range("N74").Value = -10762.77

$ws.Range("H77").Value = 3460.762
$ws.Range("I77").Value = 971.0345
$ws.Range("J77").Value = 9014.77
$ws.Range("K77").Value = 4855.1725
$ws.Range("L77").Value = 45073.85000000001
$ws.Range("M77").Value = -487.1724999999997
$ws.Range("N77").Value = -53809.85000000001

$ws.Range("H97").Value = 1592.7106
$ws.Range("I97").Value = 557.1786
$ws.Range("K97").Value = 557.1786
$ws.Range("M97").Value = -61.17859999999996

$ws.Range("H102").Value = 1508.2
$ws.Range("I102").Value = 1210
$ws.Range("K102").Value = 1210
$ws.Range("M102").Value = 412

$ws.Range("H110").Value = 1710.96
$ws.Range("I110").Value = 1683.7
$ws.Range("J110").Value = 1820
$ws.Range("K110").Value = 1683.7
$ws.Range("L110").Value = 1820
$ws.Range("M110").Value = 361.3
$ws.Range("N110").Value = -5910

$ws.Range("H122").Value = 1824.1818
$ws.Range("I122").Value = 1728.9412
$ws.Range("J122").Value = 2148
$ws.Range("K122").Value = 5186.8236
$ws.Range("L122").Value = 6444
$ws.Range("M122").Value = -2736.8236
$ws.Range("N122").Value = -11344

$ws.Range("H132").Value = 8103.5586
$ws.Range("I132").Value = 5944.84
$ws.Range("J132").Value = 14100
$ws.Range("K132").Value = 17834.52
$ws.Range("L132").Value = 42300
$ws.Range("M132").Value = -15304.52
$ws.Range("N132").Value = -47360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 695.5278
$ws.Range("I94").Value = 415.34616
$ws.Range("J94").Value = 1424
$ws.Range("K94").Value = 415.34616
$ws.Range("L94").Value = 1424
$ws.Range("M94").Value = 35.65384
$ws.Range("N94").Value = -2326

$ws.Range("H99").Value = 1067.1765
$ws.Range("I99").Value = 942
$ws.Range("J99").Value = 1246
$ws.Range("K99").Value = 942
$ws.Range("L99").Value = 1246
$ws.Range("M99").Value = 556
$ws.Range("N99").Value = -4242

$ws.Range("H105").Value = 3342.0833
$ws.Range("I105").Value = 3850
$ws.Range("K105").Value = 3850
$ws.Range("M105").Value = -2103

$ws.Range("H107").Value = 2173.6
$ws.Range("I107").Value = 1787.8
$ws.Range("J107").Value = 2945.2
$ws.Range("K107").Value = 1787.8
$ws.Range("L107").Value = 2945.2
$ws.Range("M107").Value = 132.2
$ws.Range("N107").Value = -6785.2

$ws.Range("H132").Value = 40702
$ws.Range("J132").Value = 40702
$ws.Range("L132").Value = 40702
$ws.Range("N132").Value = -50822

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3369.3635
$ws.Range("I16").Value = 1229.2222
$ws.Range("J16").Value = 13000
$ws.Range("K16").Value = 1229.2222
$ws.Range("L16").Value = 13000
$ws.Range("M16").Value = -942.2221999999999
$ws.Range("N16").Value = -13574

$ws.Range("H31").Value = 7096.067
$ws.Range("I31").Value = 5541.4585
$ws.Range("J31").Value = 8872.762000000001
$ws.Range("K31").Value = 5541.4585
$ws.Range("L31").Value = 8872.762000000001
$ws.Range("M31").Value = -5246.4585
$ws.Range("N31").Value = -9462.762000000001

$ws.Range("H34").Value = 7096.067
$ws.Range("I34").Value = 5541.4585
$ws.Range("J34").Value = 8872.762000000001
$ws.Range("K34").Value = 5541.4585
$ws.Range("L34").Value = 8872.762000000001
$ws.Range("M34").Value = -5339.4585
$ws.Range("N34").Value = -9276.762000000001

$ws.Range("H58").Value = 2698.238
$ws.Range("I58").Value = 2270.4285
$ws.Range("J58").Value = 2912.1428
$ws.Range("K58").Value = 2270.4285
$ws.Range("L58").Value = 2912.1428
$ws.Range("M58").Value = -2067.4285
$ws.Range("N58").Value = -3318.1428

$ws.Range("H113").Value = 3369.3635
$ws.Range("I113").Value = 1229.2222
$ws.Range("J113").Value = 13000
$ws.Range("K113").Value = 1229.2222
$ws.Range("L113").Value = 13000
$ws.Range("M113").Value = 940.7778000000001
$ws.Range("N113").Value = -17340

$ws.Range("H122").Value = 1083.3125
$ws.Range("I122").Value = 951.9
$ws.Range("J122").Value = 1302.3334
$ws.Range("K122").Value = 2855.7
$ws.Range("L122").Value = 3907.0002
$ws.Range("M122").Value = -405.6999999999998
$ws.Range("N122").Value = -8807.0002

$ws.Range("H136").Value = 2698.238
$ws.Range("I136").Value = 2270.4285
$ws.Range("J136").Value = 2912.1428
$ws.Range("K136").Value = 6811.2855
$ws.Range("L136").Value = 8736.428400000001
$ws.Range("M136").Value = -4261.2855
$ws.Range("N136").Value = -13836.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 111112110
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 111112110
$ws.Range("K36").Value = 0
$ws.Range("L36").ClearContents()
$ws.Range("M36").Value = 333336330
$ws.Range("N36").Value = -333336668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 48000
$ws.Range("J63").Value = 48000
$ws.Range("L63").Value = 48000
$ws.Range("N63").Value = -49372

$ws.Range("H66").Value = 48000
$ws.Range("J66").Value = 48000
$ws.Range("L66").Value = 144000
$ws.Range("N66").Value = -150864

$ws.Range("H70").Value = 4463.04
$ws.Range("J70").Value = 4386.3335
$ws.Range("L70").Value = 4386.3335
$ws.Range("N70").Value = -4926.3335

$ws.Range("H73").Value = 4463.04
$ws.Range("J73").Value = 4386.3335
$ws.Range("L73").Value = 4386.3335
$ws.Range("N73").Value = -6258.3335

$ws.Range("H80").Value = 4075.8462
$ws.Range("I80").Value = 2264
$ws.Range("J80").Value = 5208.25
$ws.Range("K80").Value = 2264
$ws.Range("L80").Value = 5208.25
$ws.Range("M80").Value = -1266
$ws.Range("N80").Value = -7204.25

$ws.Range("H83").Value = 4075.8462
$ws.Range("I83").Value = 2264
$ws.Range("J83").Value = 5208.25
$ws.Range("K83").Value = 11320
$ws.Range("L83").Value = 26041.25
$ws.Range("M83").Value = -6328
$ws.Range("N83").Value = -36025.25

$ws.Range("H97").Value = 1414.2069
$ws.Range("I97").Value = 1324.381
$ws.Range("J97").Value = 1650
$ws.Range("K97").Value = 1324.381
$ws.Range("L97").Value = 1650
$ws.Range("M97").Value = -828.3810000000001
$ws.Range("N97").Value = -2642

$ws.Range("H132").Value = 6601.6562
$ws.Range("I132").Value = 8003.864
$ws.Range("J132").Value = 3516.8
$ws.Range("K132").Value = 24011.592
$ws.Range("L132").Value = 10550.4
$ws.Range("M132").Value = -21481.592
$ws.Range("N132").Value = -15610.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 188.33333
$ws.Range("I55").Value = 161.42857
$ws.Range("J55").Value = 211.875
$ws.Range("K55").Value = 161.42857
$ws.Range("L55").Value = 211.875
$ws.Range("M55").Value = 11.57142999999999
$ws.Range("N55").Value = -557.875

$ws.Range("H61").Value = 2279.9167
$ws.Range("I61").Value = 2029.1428
$ws.Range("J61").Value = 2631
$ws.Range("K61").Value = 2029.1428
$ws.Range("L61").Value = 2631
$ws.Range("M61").Value = -1827.1428
$ws.Range("N61").Value = -3035

$ws.Range("H113").Value = 2279.9167
$ws.Range("I113").Value = 2029.1428
$ws.Range("J113").Value = 2631
$ws.Range("K113").Value = 2029.1428
$ws.Range("L113").Value = 2631
$ws.Range("M113").Value = 140.8571999999999
$ws.Range("N113").Value = -6971

$ws.Range("H133").Value = 46800
$ws.Range("J133").Value = 46800
$ws.Range("L133").Value = 46800
$ws.Range("N133").Value = -51860

$ws.Range("H136").Value = 4848.2974
$ws.Range("I136").Value = 1870.5
$ws.Range("J136").Value = 10345.77
$ws.Range("K136").Value = 5611.5
$ws.Range("L136").Value = 31037.31
$ws.Range("M136").Value = -3061.5
$ws.Range("N136").Value = -36137.31

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5744588
$ws.Range("I122").Value = 9137429
$ws.Range("J122").Value = 2856.5386
$ws.Range("K122").Value = 27412287
$ws.Range("L122").Value = 8569.6158
$ws.Range("M122").Value = -27409837
$ws.Range("N122").Value = -13469.6158

$ws.Range("H136").Value = 4368.647
$ws.Range("I136").Value = 6037.8423
$ws.Range("K136").Value = 18113.5269
$ws.Range("M136").Value = -15563.5269
